$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 14987.5
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 14987.5
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 14987.5
$ws.Range("N7").Value = -15211.5

$ws.Range("H14").Value = 14987.5
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 14987.5
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 14987.5
$ws.Range("N14").Value = -15369.5

$ws.Range("H15").Value = 115007.73
$ws.Range("I15").Value = 115007.73
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 345023.19
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -344854.19

$ws.Range("H32").Value = 642.8570999999999
$ws.Range("I32").Value = 700
$ws.Range("J32").Value = 620
$ws.Range("K32").Value = 700
$ws.Range("L32").Value = 620
$ws.Range("M32").Value = -374
$ws.Range("N32").Value = -1272

$ws.Range("H132").Value = 2382708.8
$ws.Range("I132").Value = 2858911.8
$ws.Range("J132").Value = 1692.9
$ws.Range("K132").Value = 8576735.399999999
$ws.Range("L132").Value = 5078.700000000001
$ws.Range("M132").Value = -8574205.399999999
$ws.Range("N132").Value = -10138.7

$ws.Range("H135").Value = 631.4400000000001
$ws.Range("I135").Value = 511.65216
$ws.Range("J135").Value = 2009
$ws.Range("K135").Value = 4604.869439999999
$ws.Range("L135").Value = 18081
$ws.Range("M135").Value = -2069.869439999999
$ws.Range("N135").Value = -23151

$ws.Range("H137").Value = 1337.2653
$ws.Range("I137").Value = 1095.2433
$ws.Range("J137").Value = 2083.5
$ws.Range("K137").Value = 3285.7299
$ws.Range("L137").Value = 6250.5
$ws.Range("M137").Value = -735.7299000000003
$ws.Range("N137").Value = -11350.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").ClearContents()

$ws.Range("H39").Value = 5800
$ws.Range("I39").Value = 3000
$ws.Range("J39").Value = 10000
$ws.Range("K39").Value = 3000
$ws.Range("L39").Value = 10000
$ws.Range("M39").Value = -2480
$ws.Range("N39").Value = -11040

$ws.Range("H61").Value = 1488.6154
$ws.Range("I61").Value = 939.38464
$ws.Range("J61").Value = 2587.077
$ws.Range("K61").Value = 939.38464
$ws.Range("L61").Value = 2587.077
$ws.Range("M61").Value = -727.38464
$ws.Range("N61").Value = -3011.077

$ws.Range("H136").Value = 1488.6154
$ws.Range("I136").Value = 939.38464
$ws.Range("J136").Value = 2587.077
$ws.Range("K136").Value = 2818.15392
$ws.Range("L136").Value = 7761.231000000001
$ws.Range("M136").Value = -268.1539199999997
$ws.Range("N136").Value = -12861.231

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4014
$ws.Range("I86").Value = 2872.4
$ws.Range("J86").Value = 5916.6665
$ws.Range("K86").Value = 2872.4
$ws.Range("L86").Value = 5916.6665
$ws.Range("M86").Value = -1749.4
$ws.Range("N86").Value = -8162.6665

$ws.Range("H89").Value = 4014
$ws.Range("I89").Value = 2872.4
$ws.Range("J89").Value = 5916.6665
$ws.Range("K89").Value = 14362
$ws.Range("L89").Value = 29583.3325
$ws.Range("M89").Value = -8746
$ws.Range("N89").Value = -40815.3325

$ws.Range("H134").Value = 21783.06
$ws.Range("I134").Value = 25862.586
$ws.Range("J134").Value = 3198.5557
$ws.Range("K134").Value = 77587.758
$ws.Range("L134").Value = 9595.667099999999
$ws.Range("M134").Value = -75052.758
$ws.Range("N134").Value = -14665.6671

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H39").Value = 5379.1665
$ws.Range("I39").Value = 1950
$ws.Range("J39").Value = 15666.667
$ws.Range("K39").Value = 1950
$ws.Range("L39").Value = 15666.667
$ws.Range("M39").Value = -1559
$ws.Range("N39").Value = -16448.667

$ws.Range("H49").Value = 5379.1665
$ws.Range("I49").Value = 1950
$ws.Range("J49").Value = 15666.667
$ws.Range("K49").Value = 1950
$ws.Range("L49").Value = 15666.667
$ws.Range("M49").Value = -1768
$ws.Range("N49").Value = -16030.667

$ws.Range("H74").Value = 29637.166
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 29637.166
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 29637.166
$ws.Range("N74").Value = -31385.166

$ws.Range("H77").Value = 29637.166
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 29637.166
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 88911.49800000001
$ws.Range("N77").Value = -97647.49800000001

$ws.Range("H99").Value = 36942.93
$ws.Range("I99").Value = 79033.30499999999
$ws.Range("J99").Value = 2744.5
$ws.Range("K99").Value = 79033.30499999999
$ws.Range("L99").Value = 2744.5
$ws.Range("M99").Value = -77535.30499999999
$ws.Range("N99").Value = -5740.5

$ws.Range("H126").Value = 36942.93
$ws.Range("I126").Value = 79033.30499999999
$ws.Range("J126").Value = 2744.5
$ws.Range("K126").Value = 237099.915
$ws.Range("L126").Value = 8233.5
$ws.Range("M126").Value = -234629.915
$ws.Range("N126").Value = -13173.5

$ws.Range("H134").Value = 937.1429000000001
$ws.Range("I134").Value = 891.06976
$ws.Range("J134").Value = 1267.3334
$ws.Range("K134").Value = 2673.20928
$ws.Range("L134").Value = 3802.0002
$ws.Range("M134").Value = -138.20928
$ws.Range("N134").Value = -8872.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 647.86206
$ws.Range("I113").Value = 540.53845
$ws.Range("J113").Value = 735.0625
$ws.Range("K113").Value = 1621.61535
$ws.Range("L113").Value = 2205.1875
$ws.Range("M113").Value = 548.38465
$ws.Range("N113").Value = -6545.1875

$ws.Range("H131").Value = 5158.875
$ws.Range("I131").Value = 7336.25
$ws.Range("J131").Value = 804.125
$ws.Range("K131").Value = 22008.75
$ws.Range("L131").Value = 2412.375
$ws.Range("M131").Value = -16968.75
$ws.Range("N131").Value = -12492.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 4211.125
$ws.Range("I13").Value = 337.8
$ws.Range("J13").Value = 10666.667
$ws.Range("K13").Value = 337.8
$ws.Range("L13").Value = 10666.667
$ws.Range("M13").Value = -198.8
$ws.Range("N13").Value = -10944.667

$ws.Range("H80").Value = 4026.4
$ws.Range("I80").Value = 2536.8572
$ws.Range("J80").Value = 7502
$ws.Range("K80").Value = 2536.8572
$ws.Range("L80").Value = 7502
$ws.Range("M80").Value = -1538.8572
$ws.Range("N80").Value = -9498

$ws.Range("H83").Value = 4026.4
$ws.Range("I83").Value = 2536.8572
$ws.Range("J83").Value = 7502
$ws.Range("K83").Value = 12684.286
$ws.Range("L83").Value = 37510
$ws.Range("M83").Value = -7692.286
$ws.Range("N83").Value = -47494

$ws.Range("H97").Value = 880.6429000000001
$ws.Range("I97").Value = 784.4545000000001
$ws.Range("J97").Value = 1233.3334
$ws.Range("K97").Value = 784.4545000000001
$ws.Range("L97").Value = 1233.3334
$ws.Range("M97").Value = -288.4545000000001
$ws.Range("N97").Value = -2225.3334

$ws.Range("H131").Value = 25750.5
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 25750.5
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 25750.5
$ws.Range("N131").Value = -35830.5

$ws.Range("H132").Value = 21410.412
$ws.Range("I132").Value = 30322.828
$ws.Range("J132").Value = 1914.5
$ws.Range("K132").Value = 90968.484
$ws.Range("L132").Value = 5743.5
$ws.Range("M132").Value = -88438.484
$ws.Range("N132").Value = -10803.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 16668350
$ws.Range("I61").Value = 1617.2
$ws.Range("J61").Value = 66668548
$ws.Range("K61").Value = 1617.2
$ws.Range("L61").Value = 66668548
$ws.Range("M61").Value = -1415.2
$ws.Range("N61").Value = -66668952

$ws.Range("H113").Value = 16668350
$ws.Range("I113").Value = 1617.2
$ws.Range("J113").Value = 66668548
$ws.Range("K113").Value = 1617.2
$ws.Range("L113").Value = 66668548
$ws.Range("M113").Value = 552.8
$ws.Range("N113").Value = -66672888

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 8833.444
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 8833.444
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 8833.444
$ws.Range("N14").Value = -9169.444

$ws.Range("H132").Value = 1245.8292
$ws.Range("I132").Value = 1243.1562
$ws.Range("J132").Value = 1255.3334
$ws.Range("K132").Value = 3729.4686
$ws.Range("L132").Value = 3766.0002
$ws.Range("M132").Value = -1199.4686
$ws.Range("N132").Value = -8826.0002
